$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "Carros"

$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Ano"
$ws.Range("C1").Value = "km"
$ws.Range("D1").Value = "Local"
$ws.Range("E1").Value = "Preco"
$ws.Range("F1").Value = "URL"

$target = $wb.Worksheets.Item("Planilha1")
$ws.Move($null, $target)
